# Append " (Changed main)" to the end of the first paragraph
# ("This is a Microsoft word document."), split across three separate
# runs - " (", "Changed main", ")" - exactly as the target OOXML does,
# instead of one merged run.
#
# A plain Range.InsertAfter() at a collapsed range gets coalesced into
# the neighbouring run when the package is serialized (same text, same
# run formatting). To keep the inserted text as independent <w:r>
# elements we briefly split the paragraph in two, type the new text into
# the fresh (run-less) paragraph that results, and then delete the
# paragraph mark again to rejoin the paragraphs. Because the paragraph
# mark that survives a merge is the *later* one, the freshly typed text
# is never absorbed into the preceding run, so it remains its own run.

$d = $word.ActiveDocument

function Insert-TextAsNewRun($pos, [string]$text) {
    $splitPoint = $d.Range($pos, $pos)
    $splitPoint.InsertParagraphAfter()

    $newRunStart = $pos + 1
    $newRun = $d.Range($newRunStart, $newRunStart)
    $newRun.InsertAfter($text)
    $newRunEnd = $newRunStart + $text.Length

    # Remove the paragraph mark inserted above (a single character sitting
    # right at $pos) to rejoin the two paragraphs into one again.
    $d.Range($pos, $newRunStart).Delete()

    # The text that used to start at $newRunStart now starts at $pos,
    # so its end moved back by the one deleted character.
    return $newRunEnd - 1
}

$rng = $d.Content
$rng.Find.Execute("This is a Microsoft word document.", $true, $false, $false,
                   $false, $false, $true, 1, $false, "", 0)

if ($rng.Find.Found) {
    $pos = $rng.End

    $pos = Insert-TextAsNewRun $pos " ("
    $pos = Insert-TextAsNewRun $pos "Changed main"
    $pos = Insert-TextAsNewRun $pos ")"
}
